$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.641.71"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "1.796.34"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "313.19"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").Value = "0.5352"
$ws.Range("E7").Value = "  -1.37%  "

$ws.Range("D8").Value = "0.3776"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "0.07526"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").Value = "42.49"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "20.96"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "6.178"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").Value = "7.422"
$ws.Range("E15").Value = "  +4.20%  "

$ws.Range("D16").Value = "1.794.80"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").Value = "90.17"
$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("D18").Value = "'0.00001064"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").Value = "0.06439"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").Value = "5.934"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D23").Value = "28.635.64"
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("D25").Value = "2.094"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "160.66"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").Value = "20.44"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").Value = "2.001.14"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "'122.80"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("D32").Value = "0.1025"
$ws.Range("E32").Value = "  -0.76%  "

$ws.Range("D33").Value = "5.658"
$ws.Range("E33").Value = "  -1.37%  "

$ws.Range("D34").Value = "3.681"
$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("E35").Value = "  +6.98%  "

$ws.Range("D36").Value = "0.06454"
$ws.Range("E36").Value = "  +6.55%  "

$ws.Range("D37").Value = "8.925"
$ws.Range("E37").Value = "  +3.36%  "

$ws.Range("D38").Value = "0.02301"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").Value = "5.039"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "11.35"
$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("D41").Value = "1.209"
$ws.Range("E41").Value = "  +5.06%  "

$ws.Range("D42").Value = "0.6249"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("D44").Value = "1.392"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").Value = "13.42"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "0.5886"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "3.657"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "126.11"
$ws.Range("E48").Value = "  +3.22%  "

$ws.Range("D49").Value = "1.961"
$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("D50").Value = "1.157"
$ws.Range("E50").Value = "  +2.11%  "

$ws.Range("D51").Value = "0.06897"
$ws.Range("E51").Value = "  +1.79%  "
